$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-26 09:41:05"
$wsZhCn.Range("G3").Value = "2016-01-26 09:41:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-26 09:41:15"
$wsDeDe.Range("G3").Value = "2016-01-26 09:42:06"
